# Fruta / hortaliza, semanal
# Insert the new weekly price record as row 21 (pushing the existing
# rows 21-48 down to 22-49) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 21; Excel copies the
# formatting of the row above (row 20) onto the freshly inserted row,
# which is what keeps column D's date number-format (style index 2).
$ws.Rows("21").Insert()

$ws.Cells.Item(21, 1).Value = 10
$ws.Cells.Item(21, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(21, 3).Value = "La Araucanía"
$ws.Cells.Item(21, 4).Value = 44721
$ws.Cells.Item(21, 5).Value = 9
$ws.Cells.Item(21, 6).Value = "Fruta"
$ws.Cells.Item(21, 7).Value = 100107
$ws.Cells.Item(21, 8).Value = "Otros"
$ws.Cells.Item(21, 9).Value = 100107001
$ws.Cells.Item(21, 10).Value = "Caqui"
$ws.Cells.Item(21, 11).Value = "Fuyu"
$ws.Cells.Item(21, 12).Value = "Primera"
$ws.Cells.Item(21, 13).Value = 65
$ws.Cells.Item(21, 14).Value = 20000
$ws.Cells.Item(21, 15).Value = 20000
$ws.Cells.Item(21, 16).Value = 20000
$ws.Cells.Item(21, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(21, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(21, 19).Value = 1333
$ws.Cells.Item(21, 20).Value = 15
